$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.005.23'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '3.540.09'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '604.28'
$ws.Range('E5').Value = '  -2.23%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '196.52'
$ws.Range('E6').Value = '  +4.75%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').Value = '  -0.95%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.207'
$ws.Range('E9').Value = '  -4.75%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.654'
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.94'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('E12').Value = '  -1.97%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.53'
$ws.Range('E13').Value = '  -1.88%  '
$ws.Range('D14').Value = '4.100.04'
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '608.04'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '12.88'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '19.20'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').Value = '70.157.56'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '3.535.60'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.996'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '18.04'
$ws.Range('E22').Value = '  +2.06%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.33'
$ws.Range('E23').Value = '  +4.06%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '102.52'
$ws.Range('E24').Value = '  -2.49%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.63'
$ws.Range('E25').Value = '  -2.13%  '
$ws.Range('E26').Value = '  +2.61%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.94'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E28').Value = '  -3.21%  '
$ws.Range('E29').Value = '  -2.75%  '
$ws.Range('B30').Value = 'dogwifhat'
$ws.Range('C30').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.36'
$ws.Range('E30').Value = '  +16.66%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.13'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '12.67'
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('E33').Value = '  -1.96%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '63.25'
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('D35').Value = '0.0₃0863'
$ws.Range('E35').Value = '  +9.95%  '
$ws.Range('D36').Value = '3.757.01'
$ws.Range('E36').Value = '  +5.73%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.06'
$ws.Range('E38').Value = '  -3.48%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.65'
$ws.Range('E39').Value = '  +1.96%  '
$ws.Range('E40').Value = '  -1.75%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '36.64'
$ws.Range('E41').Value = '  -2.06%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '496.07'
$ws.Range('E42').Value = '  -8.27%  '
$ws.Range('E43').Value = '  -3.79%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0458'
$ws.Range('E44').Value = '  -2.84%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.141'
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.84'
$ws.Range('E46').Value = '  -4.14%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.33'
$ws.Range('E47').Value = '  -2.30%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.01'
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('E49').Value = '  -4.28%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.000253'
$ws.Range('E50').Value = '  +3.89%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '130.57'
$ws.Range('E51').Value = '  -2.22%  '
